# Insert a new data row at row 112 (pushes existing rows 112..166 down to 113..167)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A112").EntireRow.Insert()

# Populate the newly inserted row with the new Cilantro price record
$ws.Range("A112").Value = 4
$ws.Range("B112").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C112").Value = "Los Lagos"
$ws.Range("D112").Value = 44466
$ws.Range("E112").Value = 10
$ws.Range("F112").Value = 100112040
$ws.Range("G112").Value = "Cilantro"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 150
$ws.Range("K112").Value = 13000
$ws.Range("L112").Value = 13000
$ws.Range("M112").Value = 13000
$ws.Range("N112").Value = "$/caja 36 atados"
$ws.Range("O112").Value = "Región Metropolitana"
$ws.Range("P112").Value = 361
$ws.Range("Q112").Value = 36
$ws.Range("R112").Value = "Hortaliza"
